# Update the Rarres2-Ccrl2 LR-pair sheet with refreshed TPM-derived NATMI values.
# Adds an "Inflammatory-Mac" target-cluster row to each existing sending-cluster
# block and appends a brand-new "MuSCs" sending-cluster block (rows 17-21),
# while refreshing every numeric column (E:T) for all rows with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Rarres2", "Ccrl2", "ECs", 2, 0.6666666666666666, 0.7062263333333333, 2.118679, 0.03412450835304152, 0.03412450835304151, 3, 1, 2.455298666666666, 7.365895999999999, 0.03113820382581564, 0.03113820382581565, 1.733996574598222, 15.605969171384, 0.001062575896552755, 0.001062575896552755),
    @("ECs", "Rarres2", "Ccrl2", "FAPs", 2, 0.6666666666666666, 0.7062263333333333, 2.118679, 0.03412450835304152, 0.03412450835304151, 2, 0.6666666666666666, 0.243558, 0.730674, 0.003088813084277056, 0.003088813084277056, 0.172007073294, 1.548063659646, 0.0001054042278953963, 0.0001054042278953963),
    @("ECs", "Rarres2", "Ccrl2", "Inflammatory-Mac", 2, 0.6666666666666666, 0.7062263333333333, 2.118679, 0.03412450835304152, 0.03412450835304151, 3, 1, 41.38566833333334, 124.157005, 0.5248548347808349, 0.5248548347808349, 29.22764879959944, 263.048839196395, 0.01791041319361283, 0.01791041319361282),
    @("ECs", "Rarres2", "Ccrl2", "MuSCs", 2, 0.6666666666666666, 0.7062263333333333, 2.118679, 0.03412450835304152, 0.03412450835304151, 3, 1, 0.2904403333333334, 0.871321, 0.003683376862192124, 0.003683376862192124, 0.2051166116621111, 1.846049504959, 0.000125693424501275, 0.000125693424501275),
    @("ECs", "Rarres2", "Ccrl2", "Resolving-Mac", 2, 0.6666666666666666, 0.7062263333333333, 2.118679, 0.03412450835304152, 0.03412450835304151, 3, 1, 34.476682, 103.430046, 0.4372347714468802, 0.4372347714468802, 24.34834071435933, 219.135066429234, 0.01492042161047926, 0.01492042161047926),
    @("FAPs", "Rarres2", "Ccrl2", "ECs", 3, 1, 18.279338, 54.838014, 0.883248602930037, 0.8832486029300368, 3, 1, 2.455298666666666, 7.365895999999999, 0.03113820382581564, 0.03113820382581565, 44.88123421894932, 403.931107970544, 0.0275027750269024, 0.0275027750269024),
    @("FAPs", "Rarres2", "Ccrl2", "FAPs", 3, 1, 18.279338, 54.838014, 0.883248602930037, 0.8832486029300368, 2, 0.6666666666666666, 0.243558, 0.730674, 0.003088813084277056, 0.003088813084277056, 4.452079004604, 40.068711041436, 0.002728189841399728, 0.002728189841399727),
    @("FAPs", "Rarres2", "Ccrl2", "Inflammatory-Mac", 3, 1, 18.279338, 54.838014, 0.883248602930037, 0.8832486029300368, 3, 1, 41.38566833333334, 124.157005, 0.5248548347808349, 0.5248548347808349, 756.5026198208967, 6808.52357838807, 0.4635772995612479, 0.4635772995612478),
    @("FAPs", "Rarres2", "Ccrl2", "MuSCs", 3, 1, 18.279338, 54.838014, 0.883248602930037, 0.8832486029300368, 3, 1, 0.2904403333333334, 0.871321, 0.003683376862192124, 0.003683376862192124, 5.309057021832666, 47.781513196494, 0.003253337467596017, 0.003253337467596017),
    @("FAPs", "Rarres2", "Ccrl2", "Resolving-Mac", 3, 1, 18.279338, 54.838014, 0.883248602930037, 0.8832486029300368, 3, 1, 34.476682, 103.430046, 0.4372347714468802, 0.4372347714468802, 630.210923396516, 5671.898310568645, 0.3861870010328909, 0.3861870010328909),
    @("Inflammatory-Mac", "Rarres2", "Ccrl2", "ECs", 1, 0.3333333333333333, 0.03164166666666667, 0.094925, 0.001528909738290919, 0.001528909738290919, 3, 1, 2.455298666666666, 7.365895999999999, 0.03113820382581564, 0.03113820382581565, 0.07768974197777777, 0.6992076777999999, 0.0000476075030621771, 0.00004760750306217709),
    @("Inflammatory-Mac", "Rarres2", "Ccrl2", "FAPs", 1, 0.3333333333333333, 0.03164166666666667, 0.094925, 0.001528909738290919, 0.001528909738290919, 2, 0.6666666666666666, 0.243558, 0.730674, 0.003088813084277056, 0.003088813084277056, 0.00770658105, 0.06935922945, 0.0000047225164043116, 0.000004722516404311599),
    @("Inflammatory-Mac", "Rarres2", "Ccrl2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.03164166666666667, 0.094925, 0.001528909738290919, 0.001528909738290919, 3, 1, 41.38566833333334, 124.157005, 0.5248548347808349, 0.5248548347808349, 1.309511522180556, 11.785603699625, 0.00080245566808549, 0.0008024556680854898),
    @("Inflammatory-Mac", "Rarres2", "Ccrl2", "MuSCs", 1, 0.3333333333333333, 0.03164166666666667, 0.094925, 0.001528909738290919, 0.001528909738290919, 3, 1, 0.2904403333333334, 0.871321, 0.003683376862192124, 0.003683376862192124, 0.00919001621388889, 0.08271014592499999, 0.000005631550754400987, 0.000005631550754400987),
    @("Inflammatory-Mac", "Rarres2", "Ccrl2", "Resolving-Mac", 1, 0.3333333333333333, 0.03164166666666667, 0.094925, 0.001528909738290919, 0.001528909738290919, 3, 1, 34.476682, 103.430046, 0.4372347714468802, 0.4372347714468802, 1.090899679616667, 9.81809711655, 0.0006684924999845395, 0.0006684924999845393),
    @("MuSCs", "Rarres2", "Ccrl2", "ECs", 3, 1, 1.678369333333333, 5.035108, 0.08109797897863065, 0.08109797897863062, 3, 1, 2.455298666666666, 7.365895999999999, 0.03113820382581564, 0.03113820382581565, 4.120897986307555, 37.088081876768, 0.002525245399298314, 0.002525245399298313),
    @("MuSCs", "Rarres2", "Ccrl2", "FAPs", 3, 1, 1.678369333333333, 5.035108, 0.08109797897863065, 0.08109797897863062, 2, 0.6666666666666666, 0.243558, 0.730674, 0.003088813084277056, 0.003088813084277056, 0.408780278088, 3.679022502792, 0.00025049649857762, 0.0002504964985776198),
    @("MuSCs", "Rarres2", "Ccrl2", "Inflammatory-Mac", 3, 1, 1.678369333333333, 5.035108, 0.08109797897863065, 0.08109797897863062, 3, 1, 41.38566833333334, 124.157005, 0.5248548347808349, 0.5248548347808349, 69.46043657017111, 625.1439291315401, 0.04256466635788882, 0.0425646663578888),
    @("MuSCs", "Rarres2", "Ccrl2", "MuSCs", 3, 1, 1.678369333333333, 5.035108, 0.08109797897863065, 0.08109797897863062, 3, 1, 0.2904403333333334, 0.871321, 0.003683376862192124, 0.003683376862192124, 0.4874661486297778, 4.387195337668, 0.0002987144193404314, 0.0002987144193404314),
    @("MuSCs", "Rarres2", "Ccrl2", "Resolving-Mac", 3, 1, 1.678369333333333, 5.035108, 0.08109797897863065, 0.08109797897863062, 3, 1, 34.476682, 103.430046, 0.4372347714468802, 0.4372347714468802, 57.86460578388534, 520.781452054968, 0.03545885630352547, 0.03545885630352545)
)

$arr = New-Object "object[,]" $data.Count,20
for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt 20; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$ws.Range("A2:T21").Value = $arr
